$d = $word.ActiveDocument

# Find the "Throw bottles" list item - the new sub-bullets ("Animations",
# "Collision") belong right underneath it as second-level (ilvl=1) bullets
# in the same numbered list (numId=1).
$throwBottles = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Throw bottles`r") {
        $throwBottles = $p
        break
    }
}

if ($throwBottles -eq $null) {
    throw "Could not find the 'Throw bottles' paragraph"
}

# Insert "Animations" right after "Throw bottles".
$throwBottles.Range.InsertParagraphAfter()
$animations = $throwBottles.Next()
$animations.Range.Text = "Animations"
$animations.Range.ListFormat.ListLevelNumber = 2

# Insert "Collision" right after "Animations".
$animations.Range.InsertParagraphAfter()
$collision = $animations.Next()
$collision.Range.Text = "Collision"
$collision.Range.ListFormat.ListLevelNumber = 2
